# Update the cached "datetimeFigureOut" date field text that appears on the
# slide master and every slide layout (footer placeholder) from 25/11/2022
# to 26/11/2022.
$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = "26/11/2022"

$layoutDateIndex = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $idx = $layoutDateIndex[$li]
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = "26/11/2022"
}

# Resize the title textbox on the first slide ("Les systèmes de vote") and
# shrink its font from 66pt to 60pt.
# Shape.Height/.Width/.Top/.Left are expressed in points (1 pt = 12700 EMU),
# so convert the target EMU value (1015663) to points.
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(3)
$titleShape.Height = 1015663 / 12700
$titleShape.TextFrame.TextRange.Font.Size = 60
